$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Maquete - Detalhes" row: expand the "Ambiente" note (Sem. 9 e 10 column)
$ws.Range("J3").Value = "Ambiente, Estético, etc."
$ws.Range("J3").WrapText = $true

# "Interface de usuário" row: drop the stray "Melhorias" tag from Sem. 5 e 6 ...
$ws.Range("H4").ClearContents()

# ... and fold it into the Sem. 9 e 10 note, which is now about documentation AND improvements
$ws.Range("J4").Value = "Melhorias e Documentação"
$ws.Range("J4").WrapText = $true

# "Indicador de trânsito a frente" row is renamed to "Placa de sinalização de trânsito"
$ws.Range("A7").Value = "Placa de sinalização de trânsito"
$ws.Range("A7").WrapText = $true

# ... its Sem. 5 e 6 milestone moves on from "Protótipo" to "Implementação" ...
$ws.Range("H7").Value = "Implementação"

# ... and the redundant Sem. 7 e 8 "Melhorias" tag is cleared
$ws.Range("I7").ClearContents()
